# "Added 'accepter' column & improved rendering on frontend"
#
# Target change (per the OOXML diff):
#   - header row gains a new "accepter" column (H1)
#   - two new accepted-challenge rows are appended (XML rows r="2" and r="3")
#   - row r="0" (a pre-existing, out-of-spec zero-indexed row baked into this
#     workbook by whatever external tool produced it) is also supposed to be
#     overwritten with a new data row.
#
# Note on row "0": genuine Excel workbooks are 1-indexed (row numbers start
# at 1) and this engine faithfully enforces that invariant on every write
# path (Range/Cells .Value/.Formula assignment all raise "Invalid row" for
# row 0; it is not reachable via Find/Replace, CurrentRegion, End(xlUp),
# named ranges, cross-workbook copy, etc. either). Row 0 only exists here
# because it was read verbatim from the original file and is preserved
# as-is by the engine as long as nothing attempts to rewrite it - there is
# no Excel-COM-legal way to author new content into it, so it is left
# untouched below and every other part of the diff is applied in full.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell: H1 = "accepter" -----------------------------------
$ws.Cells.Item(1, 8).Value = "accepter"

# --- New row 2: O8X0smRx, accepted by "thorn" ----------------------------
$ws.Cells.Item(2, 1).Value = "O8X0smRx"
$ws.Cells.Item(2, 2).Value = "trashboatsr"
$ws.Cells.Item(2, 3).Value = 1818
$ws.Cells.Item(2, 4).Value = 120
$ws.Cells.Item(2, 5).Value = "https://lichess.org/O8X0smRx"
$ws.Cells.Item(2, 6).Value = 2883
$ws.Cells.Item(2, 7).Value = $true
$ws.Cells.Item(2, 8).Value = "thorn"

# --- New row 3: NAFalFij, accepted by "thorn" ----------------------------
$ws.Cells.Item(3, 1).Value = "NAFalFij"
$ws.Cells.Item(3, 2).Value = "trashboatsr"
$ws.Cells.Item(3, 3).Value = 1818
$ws.Cells.Item(3, 4).Value = 130
$ws.Cells.Item(3, 5).Value = "https://lichess.org/NAFalFij"
$ws.Cells.Item(3, 6).Value = 2884
$ws.Cells.Item(3, 7).Value = $true
$ws.Cells.Item(3, 8).Value = "thorn"
